# Time Improved for Algo
# Re-run of the clustering/matching algorithm produced a faster (and slightly
# different) result set: several CV keyword rows are no longer present, the
# "must have" / "good to have" cluster matches and their scores changed, and
# the overall SCORE changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the CV-keyword rows (column E) that are no longer produced by the
#    algorithm. Search from the bottom up so row numbers of not-yet-processed
#    rows stay valid while we delete.
# ---------------------------------------------------------------------------
$keywordsToRemove = @(
    "tv : 2",
    "international : 2",
    "telecommunications : 4",
    "project : 13",
    "development : 20",
    "windows : 6",
    "sales : 2",
    "training : 5",
    "regulations : 1",
    "engineering : 1",
    "software : 9",
    "ip : 3",
    "communication : 1",
    "it : 2"
)

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = $lastRow; $r -ge 2; $r--) {
    $cellValue = $ws.Cells.Item($r, 5).Value2
    if ($keywordsToRemove -contains $cellValue) {
        $ws.Rows.Item($r).Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Update the overall SCORE for the single data row. The score columns in
#    this report are authored as text (not numbers). Briefly format as Text
#    while writing so e.g. "27.15" / "0.0" aren't coerced to numeric values,
#    then restore the original (General/Normal) cell style.
# ---------------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.15"

# ---------------------------------------------------------------------------
# 3. Recomputed "CLUSTER MUST HAVE MATCH" list (column F, rows 2-5) and its
#    score (G2). The list shrank from 5 to 4 entries and was reordered, so F6
#    is cleared.
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "git : 1"
$ws.Range("F3").Value = "analysis : 1"
$ws.Range("F4").Value = "sql : 1"
$ws.Range("F5").Value = "databases : 2"
$ws.Range("F6").ClearContents()
Set-TextValue $ws.Range("G2") "28.57"

# ---------------------------------------------------------------------------
# 4. Recomputed "CLUSTER GOOD TO HAVE MATCH" list (column H). It shrank from
#    2 entries to 1, so H3 is cleared, and its score (I2) changed.
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = "process : 1"
$ws.Range("H3").ClearContents()
Set-TextValue $ws.Range("I2") "14.28"

# ---------------------------------------------------------------------------
# 5. CLUSTER SOFT SCORE (K2) stays "0.0" - reassert as text for completeness.
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("K2") "0.0"
